$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-03-18 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-03-19 Tuesday", 2) | Out-Null
$d.Content.Find.Execute("4+11=15", $true, $false, $false, $false, $false, $true, 1, $false, "41+29=70", 2) | Out-Null
$d.Content.Find.Execute("94-71=23", $true, $false, $false, $false, $false, $true, 1, $false, "84-22=62", 2) | Out-Null
$d.Content.Find.Execute("92-74=18", $true, $false, $false, $false, $false, $true, 1, $false, "34+28=62", 2) | Out-Null
$d.Content.Find.Execute("66+23=89", $true, $false, $false, $false, $false, $true, 1, $false, "58+2=60", 2) | Out-Null
$d.Content.Find.Execute("32+19=51", $true, $false, $false, $false, $false, $true, 1, $false, "0+76=76", 2) | Out-Null
$d.Content.Find.Execute("40+16=56", $true, $false, $false, $false, $false, $true, 1, $false, "6+49=55", 2) | Out-Null
$d.Content.Find.Execute("17+32=49", $true, $false, $false, $false, $false, $true, 1, $false, "27+0=27", 2) | Out-Null
$d.Content.Find.Execute("1+3=4", $true, $false, $false, $false, $false, $true, 1, $false, "41-6=35", 2) | Out-Null
$d.Content.Find.Execute("51-44=7", $true, $false, $false, $false, $false, $true, 1, $false, "72+25=97", 2) | Out-Null
$d.Content.Find.Execute("1+42=43", $true, $false, $false, $false, $false, $true, 1, $false, "96-87=9", 2) | Out-Null
$d.Content.Find.Execute("80-27=53", $true, $false, $false, $false, $false, $true, 1, $false, "0+2=2", 2) | Out-Null
$d.Content.Find.Execute("86-44=42", $true, $false, $false, $false, $false, $true, 1, $false, "15+21=36", 2) | Out-Null
$d.Content.Find.Execute("54-34=20", $true, $false, $false, $false, $false, $true, 1, $false, "75+8=83", 2) | Out-Null
$d.Content.Find.Execute("26-25=1", $true, $false, $false, $false, $false, $true, 1, $false, "65-2=63", 2) | Out-Null
$d.Content.Find.Execute("9+82=91", $true, $false, $false, $false, $false, $true, 1, $false, "98-4=94", 2) | Out-Null
$d.Content.Find.Execute("5+15=20", $true, $false, $false, $false, $false, $true, 1, $false, "59-49=10", 2) | Out-Null
$d.Content.Find.Execute("43-37=6", $true, $false, $false, $false, $false, $true, 1, $false, "82-10=72", 2) | Out-Null
$d.Content.Find.Execute("58+39=97", $true, $false, $false, $false, $false, $true, 1, $false, "78-14=64", 2) | Out-Null
$d.Content.Find.Execute("99-35=64", $true, $false, $false, $false, $false, $true, 1, $false, "24+12=36", 2) | Out-Null
$d.Content.Find.Execute("58+40=98", $true, $false, $false, $false, $false, $true, 1, $false, "33-0=33", 2) | Out-Null
$d.Content.Find.Execute("1+36=37", $true, $false, $false, $false, $false, $true, 1, $false, "3+66=69", 2) | Out-Null
$d.Content.Find.Execute("80+7=87", $true, $false, $false, $false, $false, $true, 1, $false, "85-61=24", 2) | Out-Null
$d.Content.Find.Execute("56+32=88", $true, $false, $false, $false, $false, $true, 1, $false, "23-6=17", 2) | Out-Null
$d.Content.Find.Execute("51+12=63", $true, $false, $false, $false, $false, $true, 1, $false, "3+48=51", 2) | Out-Null
$d.Content.Find.Execute("53+11=64", $true, $false, $false, $false, $false, $true, 1, $false, "21-20=1", 2) | Out-Null
$d.Content.Find.Execute("61-28=33", $true, $false, $false, $false, $false, $true, 1, $false, "25+57=82", 2) | Out-Null
$d.Content.Find.Execute("32-10=22", $true, $false, $false, $false, $false, $true, 1, $false, "19+55=74", 2) | Out-Null
$d.Content.Find.Execute("0+75=75", $true, $false, $false, $false, $false, $true, 1, $false, "35+14=49", 2) | Out-Null
$d.Content.Find.Execute("74+23=97", $true, $false, $false, $false, $false, $true, 1, $false, "69-65=4", 2) | Out-Null
$d.Content.Find.Execute("70-56=14", $true, $false, $false, $false, $false, $true, 1, $false, "59+24=83", 2) | Out-Null
$d.Content.Find.Execute("42-39=3", $true, $false, $false, $false, $false, $true, 1, $false, "80-63=17", 2) | Out-Null
$d.Content.Find.Execute("15+84=99", $true, $false, $false, $false, $false, $true, 1, $false, "49-41=8", 2) | Out-Null
$d.Content.Find.Execute("28+37=65", $true, $false, $false, $false, $false, $true, 1, $false, "46+7=53", 2) | Out-Null
$d.Content.Find.Execute("0+68=68", $true, $false, $false, $false, $false, $true, 1, $false, "17+17=34", 2) | Out-Null
$d.Content.Find.Execute("56-56=0", $true, $false, $false, $false, $false, $true, 1, $false, "38+44=82", 2) | Out-Null
$d.Content.Find.Execute("45-10=35", $true, $false, $false, $false, $false, $true, 1, $false, "12+0=12", 2) | Out-Null
$d.Content.Find.Execute("59-34=25", $true, $false, $false, $false, $false, $true, 1, $false, "83-11=72", 2) | Out-Null
$d.Content.Find.Execute("38+15=53", $true, $false, $false, $false, $false, $true, 1, $false, "62+37=99", 2) | Out-Null
$d.Content.Find.Execute("29-28=1", $true, $false, $false, $false, $false, $true, 1, $false, "58+32=90", 2) | Out-Null
$d.Content.Find.Execute("19+74=93", $true, $false, $false, $false, $false, $true, 1, $false, "53-25=28", 2) | Out-Null
$d.Content.Find.Execute("94-72=22", $true, $false, $false, $false, $false, $true, 1, $false, "33-20=13", 2) | Out-Null
$d.Content.Find.Execute("28-15=13", $true, $false, $false, $false, $false, $true, 1, $false, "1+48=49", 2) | Out-Null
$d.Content.Find.Execute("12+1=13", $true, $false, $false, $false, $false, $true, 1, $false, "51-17=34", 2) | Out-Null
$d.Content.Find.Execute("15-11=4", $true, $false, $false, $false, $false, $true, 1, $false, "2+26=28", 2) | Out-Null
$d.Content.Find.Execute("59-53=6", $true, $false, $false, $false, $false, $true, 1, $false, "56+38=94", 2) | Out-Null
$d.Content.Find.Execute("74-38=36", $true, $false, $false, $false, $false, $true, 1, $false, "12+26=38", 2) | Out-Null
$d.Content.Find.Execute("3-1=2", $true, $false, $false, $false, $false, $true, 1, $false, "64-51=13", 2) | Out-Null
$d.Content.Find.Execute("66-16=50", $true, $false, $false, $false, $false, $true, 1, $false, "74-44=30", 2) | Out-Null
$d.Content.Find.Execute("83-24=59", $true, $false, $false, $false, $false, $true, 1, $false, "48+18=66", 2) | Out-Null
$d.Content.Find.Execute("72-3=69", $true, $false, $false, $false, $false, $true, 1, $false, "98-49=49", 2) | Out-Null
$d.Content.Find.Execute("68-16=52", $true, $false, $false, $false, $false, $true, 1, $false, "72+21=93", 2) | Out-Null
$d.Content.Find.Execute("10+42=52", $true, $false, $false, $false, $false, $true, 1, $false, "63-23=40", 2) | Out-Null
$d.Content.Find.Execute("33+29=62", $true, $false, $false, $false, $false, $true, 1, $false, "5+10=15", 2) | Out-Null
$d.Content.Find.Execute("47+4=51", $true, $false, $false, $false, $false, $true, 1, $false, "43+12=55", 2) | Out-Null
$d.Content.Find.Execute("55-14=41", $true, $false, $false, $false, $false, $true, 1, $false, "31-26=5", 2) | Out-Null
$d.Content.Find.Execute("37-28=9", $true, $false, $false, $false, $false, $true, 1, $false, "93-56=37", 2) | Out-Null
$d.Content.Find.Execute("34+45=79", $true, $false, $false, $false, $false, $true, 1, $false, "1+22=23", 2) | Out-Null
$d.Content.Find.Execute("83-21=62", $true, $false, $false, $false, $false, $true, 1, $false, "59-41=18", 2) | Out-Null
$d.Content.Find.Execute("18+42=60", $true, $false, $false, $false, $false, $true, 1, $false, "21+42=63", 2) | Out-Null
$d.Content.Find.Execute("88-23=65", $true, $false, $false, $false, $false, $true, 1, $false, "16+83=99", 2) | Out-Null
$d.Content.Find.Execute("13+80=93", $true, $false, $false, $false, $false, $true, 1, $false, "74+16=90", 2) | Out-Null
$d.Content.Find.Execute("74+24=98", $true, $false, $false, $false, $false, $true, 1, $false, "2+38=40", 2) | Out-Null
$d.Content.Find.Execute("26+35=61", $true, $false, $false, $false, $false, $true, 1, $false, "7+23=30", 2) | Out-Null
$d.Content.Find.Execute("72+7=79", $true, $false, $false, $false, $false, $true, 1, $false, "5+91=96", 2) | Out-Null
$d.Content.Find.Execute("58+10=68", $true, $false, $false, $false, $false, $true, 1, $false, "45+12=57", 2) | Out-Null
$d.Content.Find.Execute("19+52=71", $true, $false, $false, $false, $false, $true, 1, $false, "33+36=69", 2) | Out-Null
$d.Content.Find.Execute("65-31=34", $true, $false, $false, $false, $false, $true, 1, $false, "36+10=46", 2) | Out-Null
$d.Content.Find.Execute("44-35=9", $true, $false, $false, $false, $false, $true, 1, $false, "6+5=11", 2) | Out-Null
$d.Content.Find.Execute("72-36=36", $true, $false, $false, $false, $false, $true, 1, $false, "31+51=82", 2) | Out-Null
$d.Content.Find.Execute("1+66=67", $true, $false, $false, $false, $false, $true, 1, $false, "54-4=50", 2) | Out-Null
$d.Content.Find.Execute("56-12=44", $true, $false, $false, $false, $false, $true, 1, $false, "49+47=96", 2) | Out-Null
$d.Content.Find.Execute("63+35=98", $true, $false, $false, $false, $false, $true, 1, $false, "61+19=80", 2) | Out-Null
$d.Content.Find.Execute("52+37=89", $true, $false, $false, $false, $false, $true, 1, $false, "86-48=38", 2) | Out-Null
$d.Content.Find.Execute("27-12=15", $true, $false, $false, $false, $false, $true, 1, $false, "31-7=24", 2) | Out-Null
$d.Content.Find.Execute("84-26=58", $true, $false, $false, $false, $false, $true, 1, $false, "87-69=18", 2) | Out-Null
$d.Content.Find.Execute("74-19=55", $true, $false, $false, $false, $false, $true, 1, $false, "14+58=72", 2) | Out-Null
$d.Content.Find.Execute("75-49=26", $true, $false, $false, $false, $false, $true, 1, $false, "32-20=12", 2) | Out-Null
$d.Content.Find.Execute("35+49=84", $true, $false, $false, $false, $false, $true, 1, $false, "2+12=14", 2) | Out-Null
$d.Content.Find.Execute("96-42=54", $true, $false, $false, $false, $false, $true, 1, $false, "48+43=91", 2) | Out-Null
$d.Content.Find.Execute("24+69=93", $true, $false, $false, $false, $false, $true, 1, $false, "92-82=10", 2) | Out-Null
$d.Content.Find.Execute("32-21=11", $true, $false, $false, $false, $false, $true, 1, $false, "24+14=38", 2) | Out-Null
$d.Content.Find.Execute("63+14=77", $true, $false, $false, $false, $false, $true, 1, $false, "13-9=4", 2) | Out-Null
$d.Content.Find.Execute("96-51=45", $true, $false, $false, $false, $false, $true, 1, $false, "91-36=55", 2) | Out-Null
$d.Content.Find.Execute("78-25=53", $true, $false, $false, $false, $false, $true, 1, $false, "31+8=39", 2) | Out-Null
$d.Content.Find.Execute("96-9=87", $true, $false, $false, $false, $false, $true, 1, $false, "78+3=81", 2) | Out-Null
$d.Content.Find.Execute("44+29=73", $true, $false, $false, $false, $false, $true, 1, $false, "57+18=75", 2) | Out-Null
$d.Content.Find.Execute("75-55=20", $true, $false, $false, $false, $false, $true, 1, $false, "33+7=40", 2) | Out-Null
$d.Content.Find.Execute("60-52=8", $true, $false, $false, $false, $false, $true, 1, $false, "21-4=17", 2) | Out-Null
$d.Content.Find.Execute("94-77=17", $true, $false, $false, $false, $false, $true, 1, $false, "25-6=19", 2) | Out-Null
$d.Content.Find.Execute("34-2=32", $true, $false, $false, $false, $false, $true, 1, $false, "42-17=25", 2) | Out-Null
$d.Content.Find.Execute("95-20=75", $true, $false, $false, $false, $false, $true, 1, $false, "54-32=22", 2) | Out-Null
$d.Content.Find.Execute("67-14=53", $true, $false, $false, $false, $false, $true, 1, $false, "41-37=4", 2) | Out-Null
$d.Content.Find.Execute("6+92=98", $true, $false, $false, $false, $false, $true, 1, $false, "43-41=2", 2) | Out-Null
$d.Content.Find.Execute("23+5=28", $true, $false, $false, $false, $false, $true, 1, $false, "0+34=34", 2) | Out-Null
$d.Content.Find.Execute("33-18=15", $true, $false, $false, $false, $false, $true, 1, $false, "80-42=38", 2) | Out-Null
$d.Content.Find.Execute("64-50=14", $true, $false, $false, $false, $false, $true, 1, $false, "16+3=19", 2) | Out-Null
$d.Content.Find.Execute("17+34=51", $true, $false, $false, $false, $false, $true, 1, $false, "13+45=58", 2) | Out-Null
$d.Content.Find.Execute("97-52=45", $true, $false, $false, $false, $false, $true, 1, $false, "98-87=11", 2) | Out-Null
$d.Content.Find.Execute("90-35=55", $true, $false, $false, $false, $false, $true, 1, $false, "92+2=94", 2) | Out-Null
$d.Content.Find.Execute("62-17=45", $true, $false, $false, $false, $false, $true, 1, $false, "81-67=14", 2) | Out-Null
